$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update URL: matchsource -> matchsync
$ws.Range("B2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/sct-abo-group-codes"

# Set Experimental flag to the text "true" (must remain a text value, not a
# boolean, so route it through a formula + paste-values so Excel's
# auto-boolean-conversion on typed "true"/"false" literals is bypassed).
$ws.Range("B7").Formula = "=""true"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Update Date to new timestamp
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"
